$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.850.47"
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").Value = "'2.354.87"
$ws.Range("E3").Value = '  -0.23%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = "'315.14"
$ws.Range("E5").Value = '  -3.65%  '
$ws.Range("D6").Value = "'108.81"
$ws.Range("E6").Value = '  +9.00%  '
$ws.Range("D7").Value = "'0.637"
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").Value = "'0.616"
$ws.Range("E9").Value = '  +0.14%  '
$ws.Range("D10").Value = "'40.67"
$ws.Range("E10").Value = '  +1.64%  '
$ws.Range("D11").Value = "'0.0928"
$ws.Range("E11").Value = '  +1.03%  '
$ws.Range("D12").Value = "'8.53"
$ws.Range("E12").Value = '  +2.01%  '
$ws.Range("D13").Value = "'1.00"
$ws.Range("E13").Value = '  -0.58%  '
$ws.Range("D14").Value = "'0.107"
$ws.Range("E14").Value = '  +1.20%  '
$ws.Range("D15").Value = "'15.87"
$ws.Range("E15").Value = '  -1.99%  '
$ws.Range("D16").Value = "'2.711.41"
$ws.Range("E16").Value = '  -0.58%  '
$ws.Range("D17").Value = "'2.360.17"
$ws.Range("E17").Value = '  -0.28%  '
$ws.Range("D18").Value = "'42.839.85"
$ws.Range("E18").Value = '  +0.49%  '
$ws.Range("D19").Value = "'7.61"
$ws.Range("E19").Value = '  -0.84%  '
$ws.Range("D20").Value = "'0.0000106"
$ws.Range("E20").Value = '  +0.75%  '
$ws.Range("D21").Value = "'76.61"
$ws.Range("E21").Value = '  +2.13%  '
$ws.Range("D22").Value = "'3.58"
$ws.Range("E22").Value = '  -2.58%  '
$ws.Range("D23").Value = "'271.70"
$ws.Range("E23").Value = '  -1.75%  '
$ws.Range("D24").Value = "'2.33"
$ws.Range("E24").Value = '  +0.82%  '
$ws.Range("D25").Value = "'9.50"
$ws.Range("E25").Value = '  -1.31%  '
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = '  +0.41%  '
$ws.Range("D27").Value = "'11.35"
$ws.Range("E27").Value = '  -0.53%  '
$ws.Range("D28").Value = "'23.46"
$ws.Range("E28").Value = '  -1.85%  '
$ws.Range("D29").Value = "'2.26"
$ws.Range("E29").Value = '  +2.02%  '
$ws.Range("D30").Value = "'36.35"
$ws.Range("E30").Value = '  +3.58%  '
$ws.Range("D31").Value = "'166.75"
$ws.Range("E31").Value = '  -4.01%  '
$ws.Range("D32").Value = "'0.0906"
$ws.Range("E32").Value = '  +0.78%  '
$ws.Range("D33").Value = "'6.14"
$ws.Range("E33").Value = '  +4.71%  '
$ws.Range("D34").Value = "'2.92"
$ws.Range("E34").Value = '  -5.89%  '
$ws.Range("D35").Value = "'0.122"
$ws.Range("E35").Value = '  +16.87%  '
$ws.Range("D36").Value = "'0.131"
$ws.Range("E36").Value = '  -0.35%  '
$ws.Range("D37").Value = "'4.67"
$ws.Range("E37").Value = '  +1.99%  '
$ws.Range("D38").Value = "'0.0359"
$ws.Range("E38").Value = '  +0.61%  '
$ws.Range("D39").Value = "'3.80"
$ws.Range("E39").Value = '  -2.42%  '
$ws.Range("D40").Value = "'2.66"
$ws.Range("E40").Value = '  -7.33%  '
$ws.Range("D41").Value = "'106.13"
$ws.Range("E41").Value = '  +18.32%  '
$ws.Range("D42").Value = "'1.50"
$ws.Range("E42").Value = '  -0.16%  '
$ws.Range("D43").Value = "'0.237"
$ws.Range("E43").Value = '  +4.65%  '
$ws.Range("D44").Value = "'71.71"
$ws.Range("E44").Value = '  +4.90%  '
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("D46").Value = "'12.42"
$ws.Range("E46").Value = '  +4.44%  '
$ws.Range("D47").Value = "'113.79"
$ws.Range("E47").Value = '  -1.16%  '
$ws.Range("D48").Value = "'79.81"
$ws.Range("E48").Value = '  +17.71%  '
$ws.Range("D49").Value = "'5.52"
$ws.Range("E49").Value = '  +1.57%  '
$ws.Range("D50").Value = "'9.08"
$ws.Range("E50").Value = '  +1.22%  '
$ws.Range("D51").Value = "'1.28"
$ws.Range("E51").Value = '  +1.91%  '
